# Add a new cell B12 with the text "to avail all details" (this introduces a
# new shared string), and update the sheet's active selection to B12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B12").Value = "to avail all details"
$ws.Range("B12").Select()
